# Apply cryptocurrency price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.226.61"
$ws.Range("E2").Value = "  +3.02%  "
$ws.Range("D3").Value = "3.110.18"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D5").Value = "'524.07"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'144.83"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "'7.42"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "3.640.96"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").Value = "'27.23"
$ws.Range("E14").Value = "  +5.79%  "
$ws.Range("D15").Value = "'0.0000168"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "59.196.65"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "'6.22"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "3.107.24"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'8.22"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").Value = "'345.44"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").Value = "'66.09"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "0.0₃0940"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "'6.88"
$ws.Range("E28").Value = "  +6.20%  "
$ws.Range("D29").Value = "'7.31"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").Value = "'21.13"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "'155.48"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").Value = "'4.67"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("D35").Value = "'6.22"
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("D36").Value = "'26.99"
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("E37").Value = "  +5.20%  "
$ws.Range("D38").Value = "'0.0692"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "3.149.73"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").Value = "'36.84"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.668"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'1.48"
$ws.Range("E44").Value = "  +5.81%  "
$ws.Range("D45").Value = "2.300.96"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").Value = "'0.0259"
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "'20.98"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "'0.978"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "'6.07"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").Value = "'0.764"
$ws.Range("E50").Value = "  +11.07%  "
$ws.Range("D51").Value = "'263.12"
$ws.Range("E51").Value = "  +11.65%  "
